$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently holds 21 data rows (rows 2-22) for Profit_thres 3000..5000.
# We need to prepend 20 new data rows (Profit_thres 1000..2900) above the
# existing data, pushing the current rows 2-22 down to rows 22-42.

# Step 1: shift the existing data rows (2-22) down by 20 rows (to 22-42),
# copying from the bottom up so we never overwrite a row before reading it.
for ($r = 22; $r -ge 2; $r--) {
    $newRow = $r + 20
    for ($c = 1; $c -le 11; $c++) {
        $srcCell = $ws.Cells.Item($r, $c)
        $dstCell = $ws.Cells.Item($newRow, $c)
        $dstCell.Value = $srcCell.Value()
    }
}

# Step 2: fill the now-vacant rows 2-21 with the new optimization results.
$newData = @(
    @(1000, 237634, 20, 20, 0, -17823),
    @(1100, 221181, 18, 18, 0, -17176),
    @(1200, 205660, 17, 17, 0, -16102),
    @(1300, 195813, 16, 16, 0, -15029),
    @(1400, 180994, 15, 15, 0, -13956),
    @(1500, 165537, 14, 14, 0, -12882),
    @(1600, 157006, 13, 13, 0, -11808),
    @(1700, 151028, 12, 12, 0, -11808),
    @(1800, 145576, 12, 12, 0, -10735),
    @(1900, 134162, 11, 11, 0, -9235),
    @(2000, 131346, 11, 11, 0, -9662),
    @(2100, 126124, 10, 10, 0, -9662),
    @(2200, 120712, 10, 10, 0, -9662),
    @(2300, 115914, 9, 9, 0, -8588),
    @(2400, 112366, 9, 9, 0, -8588),
    @(2500, 107690, 9, 9, 0, -8588),
    @(2600, 103545, 8, 8, 0, -7514),
    @(2700, 101296, 8, 8, 0, -7514),
    @(2800, 97749, 8, 8, 0, -7514),
    @(2900, 95012, 8, 8, 0, -7514)
)

for ($idx = 0; $idx -lt $newData.Count; $idx++) {
    $row = $idx + 2
    $vals = $newData[$idx]

    $ws.Cells.Item($row, 1).Value = $false
    $ws.Cells.Item($row, 2).Value = $false
    $ws.Cells.Item($row, 3).Value = $true
    $ws.Cells.Item($row, 6).Value = $vals[0]
    $ws.Cells.Item($row, 7).Value = $vals[1]
    $ws.Cells.Item($row, 8).Value = $vals[2]
    $ws.Cells.Item($row, 9).Value = $vals[3]
    $ws.Cells.Item($row, 10).Value = $vals[4]
    $ws.Cells.Item($row, 11).Value = $vals[5]
}
